$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them to numeric values
# and the "trailing zero" formatting used by the source feed (e.g. "60.60",
# "65.00") would be lost.

$ws.Range('D2').Value = '43.877.69'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '2.358.39'
$ws.Range('E3').Value = '  +0.03%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.16'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('E6').Value = '  -1.69%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '73.32'
$ws.Range('E7').Value = '  -1.49%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.607'
$ws.Range('E9').Value = '  +2.19%  '
$ws.Range('E10').Value = '  +1.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '60.60'
$ws.Range('E11').Value = '  +5.98%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '33.83'
$ws.Range('E12').Value = '  +5.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.109'
$ws.Range('E13').Value = '  +0.88%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.22'
$ws.Range('E14').Value = '  -0.79%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.20'
$ws.Range('E15').Value = '  -2.71%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.909'
$ws.Range('E16').Value = '  +0.66%  '
$ws.Range('D17').Value = '2.363.15'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').Value = '43.914.26'
$ws.Range('E18').Value = '  -0.05%  '
$ws.Range('E19').Value = '  +0.82%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '77.75'
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('E21').Value = '  -3.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '252.86'
$ws.Range('E22').Value = '  -1.37%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('E24').Value = '  +2.65%  '
$ws.Range('E25').Value = '  -5.33%  '
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.46'
$ws.Range('E27').Value = '  -2.39%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.30'
$ws.Range('E28').Value = '  +1.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '176.36'
$ws.Range('E29').Value = '  +0.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.25'
$ws.Range('E30').Value = '  -2.21%  '
$ws.Range('E31').Value = '  +0.62%  '
$ws.Range('E32').Value = '  -2.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0747'
$ws.Range('E33').Value = '  -1.98%  '
$ws.Range('E34').Value = '  -3.48%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.33'
$ws.Range('E35').Value = '  -1.51%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.78'
$ws.Range('E36').Value = '  +0.82%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.58'
$ws.Range('E37').Value = '  +4.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.40'
$ws.Range('E38').Value = '  +1.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0275'
$ws.Range('E39').Value = '  -1.52%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.44'
$ws.Range('E40').Value = '  +14.85%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '65.00'
$ws.Range('E41').Value = '  +12.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '19.80'
$ws.Range('E42').Value = '  +3.15%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.06'
$ws.Range('E43').Value = '  -0.97%  '
$ws.Range('B44').Value = 'Cronos'
$ws.Range('C44').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.106'
$ws.Range('E44').Value = '  -7.08%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.201'
$ws.Range('E45').Value = '  -2.24%  '
$ws.Range('E46').Value = '  +0.00%  '
$ws.Range('E47').Value = '  -0.74%  '
$ws.Range('E48').Value = '  -1.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.15'
$ws.Range('E49').Value = '  -1.87%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '98.13'
$ws.Range('E50').Value = '  -2.17%  '
$ws.Range('B51').Value = 'TerraClassic'
$ws.Range('C51').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000211'
$ws.Range('E51').Value = '  +13.44%  '
